$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6170411985018727
$ws1.Range("C2").Value = 0.5669882100750268
$ws1.Range("D2").Value = 0.9906367041198502
$ws1.Range("E2").Value = 0.7211997273346966
$ws1.Range("F2").Value = 0.8618442489410232
$ws1.Range("G2").Value = 0.9629629629629629
$ws1.Range("H2").Value = 0.7298461193171457
$ws1.Range("I2").Value = 529
$ws1.Range("J2").Value = 404
$ws1.Range("K2").Value = 130
$ws1.Range("L2").Value = 5

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.9629629629629629
$ws2.Range("C2").Value = 0.2434456928838951
$ws2.Range("D2").Value = 0.3886397608370702

$ws2.Range("B3").Value = 0.5669882100750268
$ws2.Range("C3").Value = 0.9906367041198502
$ws2.Range("D3").Value = 0.7211997273346966

$ws2.Range("B4").Value = 0.6170411985018727
$ws2.Range("C4").Value = 0.6170411985018727
$ws2.Range("D4").Value = 0.6170411985018727
$ws2.Range("E4").Value = 0.6170411985018727

$ws2.Range("B5").Value = 0.7649755865189949
$ws2.Range("C5").Value = 0.6170411985018727
$ws2.Range("D5").Value = 0.5549197440858834

$ws2.Range("B6").Value = 0.7649755865189948
$ws2.Range("C6").Value = 0.6170411985018727
$ws2.Range("D6").Value = 0.5549197440858834

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 130
$ws3.Range("C2").Value = 404
$ws3.Range("B3").Value = 5
$ws3.Range("C3").Value = 529
